$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# The last row of the table currently reads:
#   05/02/2019 | Main Work | Built first draft of ER diagram to describe the database.
# Two new rows are inserted immediately above it, and the (originally last)
# row's entry text is changed to describe a different piece of work.

$origLastIndex = $table.Rows.Count
$lastRow = $table.Rows.Item($origLastIndex)

# Rows.Add(beforeRow) inserts a new (blank) row immediately above beforeRow.
# Insert twice so the table grows by two rows directly above the original
# last row, then address the new rows (and the now-shifted original last
# row) by their final numeric index, which is stable once both inserts
# have happened.
$table.Rows.Add($lastRow) | Out-Null
$table.Rows.Add($lastRow) | Out-Null

$row1 = $table.Rows.Item($origLastIndex)
$row1.Cells.Item(1).Range.Text = "05/02/2019"
$row1.Cells.Item(2).Range.Text = "Main Work"
$row1.Cells.Item(3).Range.Text = "Built first draft of ER diagram to describe the database."

$row2 = $table.Rows.Item($origLastIndex + 1)
$row2.Cells.Item(1).Range.Text = "05/02/2019"
$row2.Cells.Item(2).Range.Text = "Idea"
$row2.Cells.Item(3).Range.Text = "Thinking of changing project title to “Live Lecture Feedback System” because traffic lights might not be the best form of feedback."

$row3 = $table.Rows.Item($origLastIndex + 2)
$row3.Cells.Item(3).Range.Text = "Reformatted project outline's bibliography."
